# Rename the workbook's two sheets to lower-case tab names:
#   "Metadata" -> "metadata"
#   "Lung"     -> "lung"
$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("Metadata").Name = "metadata"
$wb.Worksheets.Item("Lung").Name = "lung"
